$d = $word.ActiveDocument

# The <id>...</id> tag pairs were previously split across several runs
# (e.g. "<id>" / "p106v_" / "1" / "</id>"). Collapse each one back into a
# single run/text-node so the visible text reads "<id>p106v_1</id>" etc.
# Using Find/Execute with Replace:=wdReplaceAll (2) on text that spans the
# runs merges them into one run that carries the formatting of the first
# (the "<id>" run), which is exactly what's wanted here.

$wdReplaceAll = 2
$wdFindContinue = 1

$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Execute("<id>p106v_1</id>", $true, $false, $false, $false, $false, $true, `
                  $wdFindContinue, $false, "<id>p106v_1</id>", $wdReplaceAll)

$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute("<id>p107r_1</id>", $true, $false, $false, $false, $false, $true, `
                  $wdFindContinue, $false, "<id>p107r_1</id>", $wdReplaceAll)
